# Regenerate s_vals data to filter save games: update B2:E9 and G2:G9
# (F column / Win flags are left untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.3464964993005633, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 4.266567469086664)
    3 = @(0.3464964993005633, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1.336873824401267)
    4 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 12.0302756157461)
    5 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 6.48142807727062, 28.30127388105354)
    6 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    7 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 12.0302756157461)
    8 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 22.31973251085698)
    9 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}
